$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the crypto price/volume refresh diff.
# Column D holds price strings that look numeric ("593.57", "1.00", ...);
# force text format first so Excel does not silently coerce them to numbers/dates.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.268.19'
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.978.39'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.57'
$ws.Range("E5").Value = '  +1.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.20'
$ws.Range("E6").Value = '  -1.84%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.513'
$ws.Range("E8").Value = '  -1.87%  '
$ws.Range("B9").Value = 'LidoStakedEther'
$ws.Range("C9").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.974.66'
$ws.Range("E9").Value = '  -1.01%  '
$ws.Range("E10").Value = '  -0.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.02'
$ws.Range("E11").Value = '  +4.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.452'
$ws.Range("E12").Value = '  +2.66%  '
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.08'
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("E15").Value = '  +2.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.474.19'
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("E17").Value = '  -1.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.233.14'
$ws.Range("E18").Value = '  -1.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.976.48'
$ws.Range("E19").Value = '  -1.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '447.39'
$ws.Range("E20").Value = '  -2.69%  '
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("E22").Value = '  +0.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.34'
$ws.Range("E23").Value = '  -1.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.65'
$ws.Range("E24").Value = '  +0.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.52'
$ws.Range("E25").Value = '  +6.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.17'
$ws.Range("E26").Value = '  -2.56%  '
$ws.Range("E27").Value = '  -1.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("E29").Value = '  +2.36%  '
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.15'
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.11'
$ws.Range("E33").Value = '  -1.88%  '
$ws.Range("E34").Value = '  +1.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0813'
$ws.Range("E35").Value = '  +3.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.76'
$ws.Range("E37").Value = '  +0.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '50.08'
$ws.Range("E38").Value = '  +0.13%  '
$ws.Range("B39").Value = 'Cosmos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.00'
$ws.Range("E39").Value = '  +0.29%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.03'
$ws.Range("E40").Value = '  -2.68%  '
$ws.Range("E41").Value = '  -1.53%  '
$ws.Range("E42").Value = '  +6.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '386.73'
$ws.Range("E43").Value = '  +2.37%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("E45").Value = '  -1.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.11'
$ws.Range("E46").Value = '  +2.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.685.11'
$ws.Range("E47").Value = '  -2.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.03'
$ws.Range("E48").Value = '  +1.99%  '
$ws.Range("E50").Value = '  -0.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.14'
$ws.Range("E51").Value = '  -0.14%  '
